$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old rows 8-10 (ECs sending-cluster rows), which are no longer present
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(8).Delete()

# Update remaining data rows (2-7) with the new TPM-derived values
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Pdgfc"
$ws.Range("C2").Value = "Pdgfra"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.441874
$ws.Range("H2").Value = 7.325622
$ws.Range("I2").Value = 0.7339587032246254
$ws.Range("J2").Value = 0.7339587032246254
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.1178306666666667
$ws.Range("N2").Value = 0.353492
$ws.Range("O2").Value = 0.0005211703885903252
$ws.Range("P2").Value = 0.0005211703885903251
$ws.Range("Q2").Value = 0.287727641336
$ws.Range("R2").Value = 2.589548772024
$ws.Range("S2").Value = 0.0003825175425688292
$ws.Range("T2").Value = 0.0003825175425688291
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Pdgfc"
$ws.Range("C3").Value = "Pdgfra"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.441874
$ws.Range("H3").Value = 7.325622
$ws.Range("I3").Value = 0.7339587032246254
$ws.Range("J3").Value = 0.7339587032246254
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 225.778076
$ws.Range("N3").Value = 677.3342279999999
$ws.Range("O3").Value = 0.9986266812609277
$ws.Range("P3").Value = 0.9986266812609277
$ws.Range("Q3").Value = 551.321613554424
$ws.Range("R3").Value = 4961.894521989816
$ws.Range("S3").Value = 0.7329507439837818
$ws.Range("T3").Value = 0.7329507439837818
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Pdgfc"
$ws.Range("C4").Value = "Pdgfra"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.441874
$ws.Range("H4").Value = 7.325622
$ws.Range("I4").Value = 0.7339587032246254
$ws.Range("J4").Value = 0.7339587032246254
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.192661
$ws.Range("N4").Value = 0.5779829999999999
$ws.Range("O4").Value = 0.0008521483504820529
$ws.Range("P4").Value = 0.0008521483504820528
$ws.Range("Q4").Value = 0.4704538867139999
$ws.Range("R4").Value = 4.234084980425999
$ws.Range("S4").Value = 0.0006254416982748111
$ws.Range("T4").Value = 0.0006254416982748111
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Pdgfc"
$ws.Range("C5").Value = "Pdgfra"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.885117
$ws.Range("H5").Value = 2.655351
$ws.Range("I5").Value = 0.2660412967753745
$ws.Range("J5").Value = 0.2660412967753745
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.1178306666666667
$ws.Range("N5").Value = 0.353492
$ws.Range("O5").Value = 0.0005211703885903252
$ws.Range("P5").Value = 0.0005211703885903251
$ws.Range("Q5").Value = 0.104293926188
$ws.Range("R5").Value = 0.938645335692
$ws.Range("S5").Value = 0.000138652846021496
$ws.Range("T5").Value = 0.0001386528460214959
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Pdgfc"
$ws.Range("C6").Value = "Pdgfra"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.885117
$ws.Range("H6").Value = 2.655351
$ws.Range("I6").Value = 0.2660412967753745
$ws.Range("J6").Value = 0.2660412967753745
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 225.778076
$ws.Range("N6").Value = 677.3342279999999
$ws.Range("O6").Value = 0.9986266812609277
$ws.Range("P6").Value = 0.9986266812609277
$ws.Range("Q6").Value = 199.840013294892
$ws.Range("R6").Value = 1798.560119654028
$ws.Range("S6").Value = 0.2656759372771458
$ws.Range("T6").Value = 0.2656759372771457
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Pdgfc"
$ws.Range("C7").Value = "Pdgfra"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.885117
$ws.Range("H7").Value = 2.655351
$ws.Range("I7").Value = 0.2660412967753745
$ws.Range("J7").Value = 0.2660412967753745
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.192661
$ws.Range("N7").Value = 0.5779829999999999
$ws.Range("O7").Value = 0.0008521483504820529
$ws.Range("P7").Value = 0.0008521483504820528
$ws.Range("Q7").Value = 0.170527526337
$ws.Range("R7").Value = 1.534747737033
$ws.Range("S7").Value = 0.0002267066522072417
$ws.Range("T7").Value = 0.0002267066522072416
